# Atualizado por script em 11-11-2023 20:45
#
# 1) Rows 2 and 3 (Khenchela/El Bayadh fixtures) had the odds/result data
#    for the two matches mixed up - swap columns F..V between them.
# 2) Rows 16/17/18 had the odds/result data cyclically shifted the same
#    way - rotate columns F..V: 16 <- 18, 17 <- 16, 18 <- 17.
# 3) Four new fixtures were scraped and appended as rows 35-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the odds/result data (columns F..V) between row 2 and row 3.
# ---------------------------------------------------------------------------
$row2 = $ws.Range("F2:V2").Value()
$row3 = $ws.Range("F3:V3").Value()
$ws.Range("F2:V2").Value = $row3
$ws.Range("F3:V3").Value = $row2

# ---------------------------------------------------------------------------
# 2) Cyclically rotate the odds/result data (columns F..V) across
#    rows 16, 17 and 18: 16 <- 18, 17 <- 16, 18 <- 17.
# ---------------------------------------------------------------------------
$row16 = $ws.Range("F16:V16").Value()
$row17 = $ws.Range("F17:V17").Value()
$row18 = $ws.Range("F18:V18").Value()
$ws.Range("F16:V16").Value = $row18
$ws.Range("F17:V17").Value = $row16
$ws.Range("F18:V18").Value = $row17

# ---------------------------------------------------------------------------
# 3) Append four new fixture rows (35-38), copying the number/date format
#    from the existing data rows so the new cells keep the same style
#    (bold/bordered index column, date-time formatted match-date column).
# ---------------------------------------------------------------------------
$ws.Range("A34").Copy()
$ws.Range("A35:A38").PasteSpecial(-4122)

$ws.Range("E34").Copy()
$ws.Range("E35:E38").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$newRows = @(
    @(34, "algeria", "ligue-1", "2023-2024", 45241.625, "Khenchela", 0, "Biskra", 1, 1.74, "10/11/2023 03:13", 1.29, "11/11/2023 10:24", 3.19, "10/11/2023 03:13", 4.87, "11/11/2023 14:54", 4.9, "10/11/2023 03:13", 13.77, "11/11/2023 14:54", "https://www.betexplorer.com/football/algeria/ligue-1/khenchela-biskra/GbL62yef/"),
    @(35, "algeria", "ligue-1", "2023-2024", 45241.625, "Magra", 3, "Ben Aknoun", 1, 1.61, "11/11/2023 10:12", 1.62, "11/11/2023 14:48", 3.51, "11/11/2023 10:12", 3.56, "11/11/2023 14:48", 5.97, "11/11/2023 10:12", 6.3, "11/11/2023 14:48", "https://www.betexplorer.com/football/algeria/ligue-1/magra-es-ben-aknoun/lCJE0FP6/"),
    @(36, "algeria", "ligue-1", "2023-2024", 45241.70833333334, "ASO Chlef", 0, "MC Alger", 1, 2.11, "19/10/2023 06:12", 2.68, "11/11/2023 16:47", 2.84, "19/10/2023 06:12", 2.8, "11/11/2023 16:20", 3.64, "19/10/2023 06:12", 3.05, "11/11/2023 16:47", "https://www.betexplorer.com/football/algeria/ligue-1/aso-chlef-mc-alger/bs15NxXJ/"),
    @(37, "algeria", "ligue-1", "2023-2024", 45241.75, "Saoura", 0, "El Bayadh", 0, 1.72, "10/11/2023 06:12", 1.55, "11/11/2023 17:14", 3.15, "10/11/2023 06:12", 3.62, "11/11/2023 17:34", 4.94, "10/11/2023 06:12", 7.47, "11/11/2023 17:34", "https://www.betexplorer.com/football/algeria/ligue-1/saoura-el-bayadh/pvDNbgfJ/")
)

$r = 35
foreach ($rowData in $newRows) {
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
    $r = $r + 1
}
